# distance_state_minimal_network_louisiana.xlsx - "Add files via upload"
#
# The sheet lists pairs of neighboring states (columns B/C) with a distance
# (column D). Row 19 had B19="Florida"/C19="Louisiana" and the re-upload
# swapped those two entries to B19="Louisiana"/C19="Florida". The cursor
# position saved with the sheet also moved from G21 to D23.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the state names in B19 and C19 (Florida <-> Louisiana)
$ws.Range("B19").Value = "Louisiana"
$ws.Range("C19").Value = "Florida"

# Move the saved selection to D23 (matches the new <selection> in the sheet)
$ws.Range("D23").Select()

# Best-effort: widen the saved workbook window (cosmetic UI chrome only)
$wb.Windows.Item(1).Width = 33020
